# Applies cryptocurrency price/volume updates described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price text (e.g. "24.70", "43.003.63") that must
# stay as literal text (matching the source inlineStr cells), not be auto-coerced
# into numbers by Excel's input parsing. Temporarily mark the price column as
# Text before writing the values, then restore the default "Normal" style so no
# extra number formatting is left behind on the cells.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# --- Column B / C / D / E updates ---
$ws.Range('D2').Value = '43.003.63'
$ws.Range('E2').Value = '  +0.06%  '
$ws.Range('D3').Value = '2.329.83'
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').Value = '303.19'
$ws.Range('E5').Value = '  -0.42%  '
$ws.Range('D6').Value = '95.81'
$ws.Range('E6').Value = '  -1.47%  '
$ws.Range('D7').Value = '0.503'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('E9').Value = '  -1.09%  '
$ws.Range('D10').Value = '34.18'
$ws.Range('E10').Value = '  -3.38%  '
$ws.Range('D11').Value = '19.17'
$ws.Range('E11').Value = '  +2.11%  '
$ws.Range('D12').Value = '0.0785'
$ws.Range('E12').Value = '  -0.45%  '
$ws.Range('E13').Value = '  +3.00%  '
$ws.Range('E14').Value = '  -2.30%  '
$ws.Range('D15').Value = '2.693.48'
$ws.Range('E15').Value = '  +1.06%  '
$ws.Range('D16').Value = '2.262.12'
$ws.Range('E16').Value = '  -1.98%  '
$ws.Range('E17').Value = '  +1.15%  '
$ws.Range('D18').Value = '42.962.60'
$ws.Range('E18').Value = '  +0.22%  '
$ws.Range('D19').Value = '12.17'
$ws.Range('E19').Value = '  -3.53%  '
$ws.Range('D20').Value = '6.18'
$ws.Range('E20').Value = '  +2.16%  '
$ws.Range('E21').Value = '  -0.55%  '
$ws.Range('D22').Value = '68.01'
$ws.Range('E22').Value = '  +0.35%  '
$ws.Range('D23').Value = '236.86'
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('E24').Value = '  +4.09%  '
$ws.Range('E25').Value = '  -0.15%  '
$ws.Range('D26').Value = '2.42'
$ws.Range('E26').Value = '  -0.26%  '
$ws.Range('D27').Value = '24.70'
$ws.Range('E27').Value = '  -0.78%  '
$ws.Range('D28').Value = '2.06'
$ws.Range('E28').Value = '  -5.90%  '
$ws.Range('D29').Value = '9.14'
$ws.Range('E29').Value = '  +0.84%  '
$ws.Range('D30').Value = '31.66'
$ws.Range('E30').Value = '  -3.52%  '
$ws.Range('E31').Value = '  +0.03%  '
$ws.Range('E32').Value = '  +0.68%  '
$ws.Range('B33').Value = 'Celestia'
$ws.Range('C33').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D33').Value = '17.81'
$ws.Range('E33').Value = '  -1.87%  '
$ws.Range('B34').Value = 'Monero'
$ws.Range('C34').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D34').Value = '136.54'
$ws.Range('E34').Value = '  -17.72%  '
$ws.Range('E35').Value = '  +1.89%  '
$ws.Range('D36').Value = '4.40'
$ws.Range('E36').Value = '  -2.07%  '
$ws.Range('E37').Value = '  -1.80%  '
$ws.Range('E38').Value = '  +2.86%  '
$ws.Range('E39').Value = '  +0.78%  '
$ws.Range('D40').Value = '22.49'
$ws.Range('E40').Value = '  +25.92%  '
$ws.Range('E41').Value = '  -0.21%  '
$ws.Range('E42').Value = '  -0.48%  '
$ws.Range('D43').Value = '1.933.39'
$ws.Range('E43').Value = '  -3.15%  '
$ws.Range('D44').Value = '0.0280'
$ws.Range('E44').Value = '  -0.11%  '
$ws.Range('E45').Value = '  -2.94%  '
$ws.Range('E46').Value = '  -2.65%  '
$ws.Range('E47').Value = '  -0.78%  '
$ws.Range('B48').Value = 'HuobiToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D48').Value = '2.88'
$ws.Range('E48').Value = '  +1.80%  '
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = '2.561.95'
$ws.Range('E49').Value = '  +1.21%  '
$ws.Range('E50').Value = '  +0.23%  '
$ws.Range('E51').Value = '  +1.79%  '

# Restore the default style on column D now that the text values are set.
$priceRange.Style = "Normal"

